# Generate Report for Handback
# Update the timestamp strings recorded for the handoff/handback xliff
# generation events, as reflected in the updated shared strings.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for first row (G2)
$overview.Range("G2").Value = "2016-09-01 21:14:31"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2) for first row
$zhcn.Range("H2").Value = "2016-09-01 21:14:26"
$zhcn.Range("K2").Value = "2016-09-01 21:14:46"

# de-de sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2) for first row
$dede.Range("H2").Value = "2016-09-01 21:14:31"
$dede.Range("K2").Value = "2016-09-01 21:14:53"
